$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 853
$ws.Range("F2").Value = 176

$ws.Range("E3").Value = 435
$ws.Range("F3").Value = 5

$ws.Range("E4").Value = 307
$ws.Range("F4").Value = 321

$ws.Range("E5").Value = 677
$ws.Range("F5").Value = 506
